$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 120 (shifts old rows 120-262 down to 121-263).
$ws.Rows.Item(120).Insert()

# Populate the new row 120 with a standard data row, pointing at the new
# weekly price observation (date 44601, volume 160).
$ws.Range("A120").Value = 3
$ws.Range("B120").Value = "Femacal de La Calera"
$ws.Range("C120").Value = "Coquimbo"
$ws.Range("D120").Value = 44601
$ws.Range("E120").Value = 5
$ws.Range("F120").Value = 100112039
$ws.Range("G120").Value = "Ciboulette"
$ws.Range("H120").Value = "Sin especificar"
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 160
$ws.Range("K120").Value = 1500
$ws.Range("L120").Value = 1500
$ws.Range("M120").Value = 1500
$ws.Range("N120").Value = "`$/docena de atados"
$ws.Range("O120").Value = "Provincia de Quillota"
$ws.Range("P120").Value = 500
$ws.Range("Q120").Value = 3
$ws.Range("R120").Value = "Hortaliza"
